$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The scraper re-ran and added two new columns (height, weight) right where the old
# "fantasy points" column (E) used to be, pushing the existing fantasy-points values
# out to the new column G.

# Capture the existing "fantasy points" values (column E) before overwriting them.
$lastRow = 16
$fantasyPoints = @{}
for ($r = 2; $r -le $lastRow; $r++) {
    $fantasyPoints[$r] = $ws.Cells.Item($r, 5).Value2
}

# Update header row: E1 = height, F1 = weight, G1 = fantasy points
$ws.Range("E1").Value = "height"
$ws.Range("F1").Value = "weight"
$ws.Range("G1").Value = "fantasy points"

# New header cells (F1, G1) need the same bold/border/centered style as the
# rest of the header row (copy the formatting that E1 already carries).
$ws.Range("E1").Copy()
$ws.Range("F1:G1").PasteSpecial(-4122)

# Fill in the new height/weight columns, and move the old fantasy points values to column G
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 5).Value = 6.5
    $ws.Cells.Item($r, 6).Value = 261
    $ws.Cells.Item($r, 7).Value = $fantasyPoints[$r]
}
